# Updated main GSC export data:
# - Append three new days (2025-12-23, 2025-12-24, 2025-12-25) to the
#   "Chart" sheet's daily video-indexing table.
# - Refresh the "Table" sheet's aggregate video count to match.

$wb = $excel.ActiveWorkbook

$chart = $wb.Worksheets.Item("Chart")
$table = $wb.Worksheets.Item("Table")

# --- Chart sheet: append rows 81-83 ---------------------------------------
$chart.Range("A81:A83").NumberFormat = "@"

$chart.Range("A81").Value = "2025-12-23"
$chart.Range("B81").Value = 22.0
$chart.Range("C81").Value = 1.0
$chart.Range("D81").Value = 0.0

$chart.Range("A82").Value = "2025-12-24"
$chart.Range("B82").Value = 22.0
$chart.Range("C82").Value = 1.0
$chart.Range("D82").Value = 0.0

$chart.Range("A83").Value = "2025-12-25"
$chart.Range("B83").Value = 22.0
$chart.Range("C83").Value = 1.0
$chart.Range("D83").Value = ""

# --- Table sheet: refresh the aggregate "Videos" count --------------------
$table.Range("C2").Value = 22.0
